$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text cells) - Volume/Number and the week range.
# We patch only the specific runs of characters that changed, using the
# known character offsets of the existing rich-text runs.
# ---------------------------------------------------------------------------

# A8 = "Volume 32   Number  42" -> "...43"  (the "42" run starts at char 21, length 2)
$ws.Range("A8").Characters(21, 2).Text = "43"

# C9 = "Report Covering the Week  10/13/2025  Through  10/19/2025"
#   "10/13/2025" run starts at char 27, length 10 -> "10/20/2025"
#   "10/19/2025" run starts at char 48, length 10 -> "10/26/2025"
$ws.Range("C9").Characters(27, 10).Text = "10/20/2025"
$ws.Range("C9").Characters(48, 10).Text = "10/26/2025"

# ---------------------------------------------------------------------------
# Column widths for columns I (9) and J (10): 7.433768 -> 6.168446
# (ColumnWidth 5.43 is the value already used by the other columns that are
# stored with width 6.168446, e.g. columns C/D/F/G.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 5.43
$ws.Columns.Item(10).ColumnWidth = 5.43

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 2
$ws.Range("C15").Copy($ws.Range("G15"))
$ws.Range("E15").Copy($ws.Range("H15"))
$ws.Range("N15").Value = 0

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -32.692307692307
$ws.Range("L16").Value = -33.962264150943
$ws.Range("M16").Value = -45.3125
$ws.Range("N16").Value = -86.590038314176

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 106
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 34.177215189873
$ws.Range("L17").Value = 29.268292682926
$ws.Range("M17").Value = 130.434782608696
$ws.Range("N17").Value = 9.278350515463

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 34
$ws.Range("G18").Value = 34
$ws.Range("I18").Value = 256
$ws.Range("J18").Value = 252
$ws.Range("K18").Value = 1.587301587301
$ws.Range("L18").Value = 1.992031872509
$ws.Range("M18").Value = 25.490196078431
$ws.Range("N18").Value = -69.560047562425

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 14.705882352941
$ws.Range("I19").Value = 410
$ws.Range("J19").Value = 384
$ws.Range("K19").Value = 6.770833333333
$ws.Range("L19").Value = -20.542635658914
$ws.Range("M19").Value = 37.123745819398
$ws.Range("N19").Value = -11.062906724511

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -62.5
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -19.230769230769
$ws.Range("I20").Value = 257
$ws.Range("J20").Value = 217
$ws.Range("K20").Value = 18.433179723502
$ws.Range("L20").Value = 62.658227848101
$ws.Range("M20").Value = 117.796610169492
$ws.Range("N20").Value = -90.449646971386

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -38.709677419354
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -2.727272727272
$ws.Range("I21").Value = 1077
$ws.Range("J21").Value = 993
$ws.Range("K21").Value = 8.459214501510
$ws.Range("L21").Value = 0.560224089635
$ws.Range("M21").Value = 46.331521739130
$ws.Range("N21").Value = -75.326460481099

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -55
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = -41.818181818181
$ws.Range("I24").Value = 447
$ws.Range("J24").Value = 467
$ws.Range("K24").Value = -4.282655246252
$ws.Range("L24").Value = -7.453416149068
$ws.Range("M24").Value = 21.798365122615

# ---------------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -35.714285714285
$ws.Range("I25").Value = 85
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = -15
$ws.Range("L25").Value = -17.475728155339

# ---------------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 171
$ws.Range("J26").Value = 181
$ws.Range("K26").Value = -5.524861878453
$ws.Range("L26").Value = -9.523809523809
$ws.Range("M26").Value = 18.75

# ---------------------------------------------------------------------------
# Row 27 (UCR Rape*)
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = 2
$ws.Range("C27").Copy($ws.Range("G27"))
$ws.Range("E27").Copy($ws.Range("H27"))

# ---------------------------------------------------------------------------
# Row 28 (Other Sex Crimes)
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50

# ---------------------------------------------------------------------------
# Row 31 (Hate Crimes)
# ---------------------------------------------------------------------------
$ws.Range("F29").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 7
$ws.Range("K31").Value = 16.666666666666
$ws.Range("L31").Value = 600

# ---------------------------------------------------------------------------
# Row 33 (Traffic Fatalities)
# ---------------------------------------------------------------------------
$ws.Range("C33").Copy($ws.Range("D33"))
$ws.Range("E22").Copy($ws.Range("E33"))
